# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
# Mirrors the commit: the original scrape only pulled team statistics, not
# the season record, so three new trailing columns (AD:AF) are appended
# with the team's 2006 record (86 wins, 76 losses, 0 ties) repeated on
# every player row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new labels in AD1:AF1.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the look of the existing header cells (bold font + thin border +
# centered alignment) by copying the formatting from an existing header
# cell onto the new header cells.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows (2-55): same record value repeated for every player.
$lastRow = 55
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 86
    $ws.Cells.Item($r, 31).Value = 76
    $ws.Cells.Item($r, 32).Value = 0
}
